# Populate the DDT (data-driven testing) report sheet with the outcome of the
# cucumber-bdd-framework Excel-utility run: which rows actually executed
# ("Execute" column) and the resulting status ("Status" column), plus two
# price corrections that were fixed while wiring up the new scenario.
#
# Writes are ordered so the shared-string table is appended to in the same
# sequence the values were typed during the real editing session:
# N, $16.40, SKIPPED, FAIL, $28.99, PASS.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Faded Short Sleeve T-shirts): this scenario was not executed.
$ws.Range("A2").Value = "N"

# Row 6 (Printed Chiffon Dress): corrected price.
$ws.Range("C6").Value = "$16.40"

# Row 2 status: skipped (not executed).
$ws.Range("G2").Value = "SKIPPED"

# Row 5 (Printed Summer Dress): this one ran but failed.
$ws.Range("G5").Value = "FAIL"

# Row 5: corrected price.
$ws.Range("C5").Value = "$28.99"

# Row 3 (Blouse): ran and passed.
$ws.Range("G3").Value = "PASS"

# Row 4 (Printed Dress): not executed -> skipped.
$ws.Range("A4").Value = "N"
$ws.Range("G4").Value = "SKIPPED"

# Row 6 (Printed Chiffon Dress): not executed -> skipped.
$ws.Range("A6").Value = "N"
$ws.Range("G6").Value = "SKIPPED"

# Leave the active cell/selection on C5, matching the last cell touched.
$ws.Range("C5").Select() | Out-Null
